$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from the (now-shifted) old first data column into the new D:E columns,
# but only for the row-ranges that actually contain data, so header-only rows stay untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

# Populate the two new columns (D = newest quarter 2018-11-30, E = prior quarter 2018-08-30) with data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4226000
$ws.Range("E8").Value = 3683000
$ws.Range("D9").Value = 2325000
$ws.Range("E9").Value = 2028000
$ws.Range("D10").Value = 1901000
$ws.Range("E10").Value = 1655000
$ws.Range("D12").Value = 289000
$ws.Range("E12").Value = 269000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 12000
$ws.Range("E14").Value = 28000
$ws.Range("D15").Value = 223000
$ws.Range("E15").Value = 188000
$ws.Range("D17").Value = 3628000
$ws.Range("E17").Value = 3193000
$ws.Range("D18").Value = 598000
$ws.Range("E18").Value = 490000
$ws.Range("D20").Value = 108000
$ws.Range("E20").Value = 65000
$ws.Range("D21").Value = 929000
$ws.Range("E21").Value = 743000
$ws.Range("D22").Value = 20000
$ws.Range("E22").Value = 22000
$ws.Range("D23").Value = 686000
$ws.Range("E23").Value = 533000
$ws.Range("D24").Value = 120000
$ws.Range("E24").Value = 94000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 566000
$ws.Range("E26").Value = 439000
$ws.Range("D27").Value = 566000
$ws.Range("E27").Value = 439000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 18000
$ws.Range("E29").Value = -3000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -108000
$ws.Range("E32").Value = -65000
$ws.Range("D33").Value = 584000
$ws.Range("E33").Value = 436000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 584000
$ws.Range("E35").Value = 436000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 7575000
$ws.Range("E41").Value = 8147000
$ws.Range("D42").Value = 1534000
$ws.Range("E42").Value = 1440000
$ws.Range("D43").Value = 22907000
$ws.Range("E43").Value = 23475000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 947000
$ws.Range("E45").Value = 928000
$ws.Range("D46").Value = 32963000
$ws.Range("E46").Value = 33990000
$ws.Range("D47").Value = 971000
$ws.Range("E47").Value = 946000
$ws.Range("D48").Value = 1724000
$ws.Range("E48").Value = 1646000
$ws.Range("D49").Value = 7109000
$ws.Range("E49").Value = 6738000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 565000
$ws.Range("E52").Value = 404000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 43332000
$ws.Range("E54").Value = 43724000
$ws.Range("D57").Value = 21843000
$ws.Range("E57").Value = 22698000
$ws.Range("D58").Value = 1998000
$ws.Range("E58").Value = 2000000
$ws.Range("D59").Value = 2063000
$ws.Range("E59").Value = 1942000
$ws.Range("D60").Value = 25904000
$ws.Range("E60").Value = 26640000
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 2042000
$ws.Range("E62").Value = 1969000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 27946000
$ws.Range("E66").Value = 28609000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 5880000
$ws.Range("E72").Value = 5296000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 15386000
$ws.Range("E76").Value = 15115000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 584000
$ws.Range("E81").Value = 436000
$ws.Range("D83").Value = 223000
$ws.Range("E83").Value = 188000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1134000
$ws.Range("E89").Value = 4670000
$ws.Range("D91").Value = -224000
$ws.Range("E91").Value = -223000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 123000
$ws.Range("E94").Value = 154000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -1769000
$ws.Range("E100").Value = 1049000
$ws.Range("D101").Value = -24000
$ws.Range("E101").Value = -26000
$ws.Range("D102").Value = -536000
$ws.Range("E102").Value = 5847000

# A handful of cells in the shifted (old) columns were also revised to reflect updated
# historical figures, not just moved - apply those corrections here.
$ws.Range("F15").Value = 180000
$ws.Range("G15").Value = 185000
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "NA"
$ws.Range("H22").Value = "NA"
$ws.Range("I22").Value = "NA"
$ws.Range("J22").Value = "NA"
$ws.Range("I91").Value = -165000
$ws.Range("J91").Value = -174000
$ws.Range("H94").Value = -862000
$ws.Range("I94").Value = -651000
$ws.Range("H102").Value = 649000
$ws.Range("I102").Value = 1806000
